# Apply bug-fix edits to "MIS Upload.xlsx":
#  - Fix header text in J1: "Final%" -> "Final %"
#  - Fix header text in L1: "AbsentDays" -> "Absent Days"
#  - Update the active selection to L2 (was C1)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J1").Value = "Final %"
$ws.Range("L1").Value = "Absent Days"

$ws.Range("L2").Select()
